# Rename the sheet (Sheet1 -> "Shreyas Iyer") and rebuild the batting log
# table: add a leading "matchNo" column and grow the data from one row
# to the full 8-match history (rows 2-9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Shreyas Iyer"

# Shift the existing teamName..result columns from B..L to C..M
$ws.Columns("A:A").Insert()

# Make room for the new rows: 4 new matches before the existing one,
# 3 new matches after it (the original row 2 ends up on row 6).
$ws.Rows("2:5").Insert()
$ws.Rows("7:9").Insert()

# Build the full 9-row x 13-column table (header + 8 matches). Numeric-looking
# values (and the blank "states" entries) are prefixed with a leading quote so
# Excel keeps storing them as text, matching the scraped source data.
$data = New-Object 'object[,]' 9,13
$data[0,0] = "matchNo"
$data[0,1] = "teamName"
$data[0,2] = "batterName"
$data[0,3] = "states"
$data[0,4] = "runs"
$data[0,5] = "balls"
$data[0,6] = "fours"
$data[0,7] = "sixes"
$data[0,8] = "sr"
$data[0,9] = "opponentTeamName"
$data[0,10] = "venue"
$data[0,11] = "date"
$data[0,12] = "result"
$data[1,0] = "Qualifier"
$data[1,1] = "Delhi Capitals"
$data[1,2] = "Shreyas Iyer"
$data[1,3] = "'"
$data[1,4] = "'30"
$data[1,5] = "'27"
$data[1,6] = "'1"
$data[1,7] = "'1"
$data[1,8] = "'111.11"
$data[1,9] = "Kolkata Knight Riders"
$data[1,10] = "Sharjah"
$data[1,11] = "October 13"
$data[1,12] = "KKR won by 3 wickets (with 1 ball remaining)"
$data[2,0] = "41st"
$data[2,1] = "Delhi Capitals"
$data[2,2] = "Shreyas Iyer"
$data[2,3] = "b Narine"
$data[2,4] = "'1"
$data[2,5] = "'5"
$data[2,6] = "'0"
$data[2,7] = "'0"
$data[2,8] = "'20.00"
$data[2,9] = "Kolkata Knight Riders"
$data[2,10] = "Sharjah"
$data[2,11] = "September 28"
$data[2,12] = "KKR won by 3 wickets (with 10 balls remaining)"
$data[3,0] = "50th"
$data[3,1] = "Delhi Capitals"
$data[3,2] = "Shreyas Iyer"
$data[3,3] = "c Gaikwad b Hazlewood"
$data[3,4] = "'2"
$data[3,5] = "'7"
$data[3,6] = "'0"
$data[3,7] = "'0"
$data[3,8] = "'28.57"
$data[3,9] = "Chennai Super Kings"
$data[3,10] = "Dubai (DSC)"
$data[3,11] = "October 04"
$data[3,12] = "Capitals won by 3 wickets (with 2 balls remaining)"
$data[4,0] = "33rd"
$data[4,1] = "Delhi Capitals"
$data[4,2] = "Shreyas Iyer"
$data[4,3] = "'"
$data[4,4] = "'47"
$data[4,5] = "'41"
$data[4,6] = "'2"
$data[4,7] = "'2"
$data[4,8] = "'114.63"
$data[4,9] = "Sunrisers Hyderabad"
$data[4,10] = "Dubai (DSC)"
$data[4,11] = "September 22"
$data[4,12] = "Capitals won by 8 wickets (with 13 balls remaining)"
$data[5,0] = "Qualifier"
$data[5,1] = "Delhi Capitals"
$data[5,2] = "Shreyas Iyer"
$data[5,3] = "c Gaikwad b Hazlewood"
$data[5,4] = "'1"
$data[5,5] = "'8"
$data[5,6] = "'0"
$data[5,7] = "'0"
$data[5,8] = "'12.50"
$data[5,9] = "Chennai Super Kings"
$data[5,10] = "Dubai (DSC)"
$data[5,11] = "October 10"
$data[5,12] = "Super Kings won by 4 wickets (with 2 balls remaining)"
$data[6,0] = "56th"
$data[6,1] = "Delhi Capitals"
$data[6,2] = "Shreyas Iyer"
$data[6,3] = "c Christian b Mohammed Siraj"
$data[6,4] = "'18"
$data[6,5] = "'18"
$data[6,6] = "'1"
$data[6,7] = "'0"
$data[6,8] = "'100.00"
$data[6,9] = "Royal Challengers Bangalore"
$data[6,10] = "Dubai (DSC)"
$data[6,11] = "October 08"
$data[6,12] = "RCB won by 7 wickets"
$data[7,0] = "36th"
$data[7,1] = "Delhi Capitals"
$data[7,2] = "Shreyas Iyer"
$data[7,3] = "st †Samson b Tewatia"
$data[7,4] = "'43"
$data[7,5] = "'32"
$data[7,6] = "'1"
$data[7,7] = "'2"
$data[7,8] = "'134.37"
$data[7,9] = "Rajasthan Royals"
$data[7,10] = "Abu Dhabi"
$data[7,11] = "September 25"
$data[7,12] = "Capitals won by 33 runs"
$data[8,0] = "46th"
$data[8,1] = "Delhi Capitals"
$data[8,2] = "Shreyas Iyer"
$data[8,3] = "'"
$data[8,4] = "'33"
$data[8,5] = "'33"
$data[8,6] = "'2"
$data[8,7] = "'0"
$data[8,8] = "'100.00"
$data[8,9] = "Mumbai Indians"
$data[8,10] = "Sharjah"
$data[8,11] = "October 02"
$data[8,12] = "Capitals won by 4 wickets (with 5 balls remaining)"

$ws.Range("A1:M9").Value = $data

# Drop back to the default "Normal" style so the text-storage hint above
# doesn't leave a lingering explicit cell format (matches the source file,
# which has no per-cell styles at all).
$ws.Range("A1:M9").Style = "Normal"

